{"js": "// Update the date line and the 25 division problems to the new values,\n// matching the author's commit (each old string is unique in the doc).\nconst replacements = [\n  [\"2026-01-15 Thursday\", \"2026-01-16 Friday\"],\n  [\"80\u00f74=\", \"50\u00f75=\"],\n  [\"50\u00f72=\", \"65\u00f74=\"],\n  [\"88\u00f75=\", \"69\u00f77=\"],\n  [\"79\u00f79=\", \"88\u00f79=\"],\n  [\"79\u00f73=\", \"14\u00f73=\"],\n  [\"87\u00f75=\", \"58\u00f77=\"],\n  [\"72\u00f78=\", \"67\u00f77=\"],\n  [\"55\u00f73=\", \"21\u00f79=\"],\n  [\"96\u00f74=\", \"56\u00f76=\"],\n  [\"58\u00f72=\", \"28\u00f76=\"],\n  [\"24\u00f72=\", \"27\u00f76=\"],\n  [\"90\u00f78=\", \"60\u00f76=\"],\n  [\"41\u00f74=\", \"92\u00f79=\"],\n  [\"58\u00f76=\", \"29\u00f75=\"],\n  [\"16\u00f72=\", \"44\u00f76=\"],\n  [\"39\u00f76=\", \"25\u00f76=\"],\n  [\"86\u00f77=\", \"64\u00f79=\"],\n  [\"45\u00f78=\", \"50\u00f74=\"],\n  [\"87\u00f77=\", \"65\u00f78=\"],\n  [\"23\u00f75=\", \"54\u00f76=\"],\n  [\"32\u00f73=\", \"37\u00f74=\"],\n  [\"26\u00f74=\", \"24\u00f73=\"],\n  [\"22\u00f73=\", \"15\u00f74=\"],\n  [\"13\u00f73=\", \"35\u00f77=\"],\n  [\"47\u00f74=\", \"42\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems to the new values,\n# matching the author's commit (each old string is unique in the doc).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-15 Thursday\", \"2026-01-16 Friday\"),\n    @(\"80\u00f74=\", \"50\u00f75=\"),\n    @(\"50\u00f72=\", \"65\u00f74=\"),\n    @(\"88\u00f75=\", \"69\u00f77=\"),\n    @(\"79\u00f79=\", \"88\u00f79=\"),\n    @(\"79\u00f73=\", \"14\u00f73=\"),\n    @(\"87\u00f75=\", \"58\u00f77=\"),\n    @(\"72\u00f78=\", \"67\u00f77=\"),\n    @(\"55\u00f73=\", \"21\u00f79=\"),\n    @(\"96\u00f74=\", \"56\u00f76=\"),\n    @(\"58\u00f72=\", \"28\u00f76=\"),\n    @(\"24\u00f72=\", \"27\u00f76=\"),\n    @(\"90\u00f78=\", \"60\u00f76=\"),\n    @(\"41\u00f74=\", \"92\u00f79=\"),\n    @(\"58\u00f76=\", \"29\u00f75=\"),\n    @(\"16\u00f72=\", \"44\u00f76=\"),\n    @(\"39\u00f76=\", \"25\u00f76=\"),\n    @(\"86\u00f77=\", \"64\u00f79=\"),\n    @(\"45\u00f78=\", \"50\u00f74=\"),\n    @(\"87\u00f77=\", \"65\u00f78=\"),\n    @(\"23\u00f75=\", \"54\u00f76=\"),\n    @(\"32\u00f73=\", \"37\u00f74=\"),\n    @(\"26\u00f74=\", \"24\u00f73=\"),\n    @(\"22\u00f73=\", \"15\u00f74=\"),\n    @(\"13\u00f73=\", \"35\u00f77=\"),\n    @(\"47\u00f74=\", \"42\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
